$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.256212005038265
$ws.Range("C2").Value = 0.4386981695876671
$ws.Range("D2").Value = 0.5658366775324168
$ws.Range("E2").Value = 0.2001110725612314
$ws.Range("G2").Value = 1.63782081524883
$ws.Range("H2").Value = 1.420924937476627
$ws.Range("I2").Value = 1.04586311862893
$ws.Range("J2").Value = 0.08388357616349573
$ws.Range("M2").Value = 0.6216688249487277
$ws.Range("N2").Value = 1.629765984340938
$ws.Range("B3").Value = 1.162872876921767
$ws.Range("C3").Value = 0.4062450462388938
$ws.Range("D3").Value = 0.5624847122200549
$ws.Range("E3").Value = 0.1998589200000573
$ws.Range("G3").Value = 1.618832522791877
$ws.Range("H3").Value = 1.418822304533904
$ws.Range("I3").Value = 1.046224515816931
$ws.Range("J3").Value = 0.08425611392980237
$ws.Range("M3").Value = 0.5956517816867688
$ws.Range("N3").Value = 1.650530304693644
$ws.Range("B4").Value = 1.106115757843043
$ws.Range("C4").Value = 0.3865349000508047
$ws.Range("D4").Value = 0.5606975081109056
$ws.Range("E4").Value = 0.1997922161827361
$ws.Range("G4").Value = 1.608302053661305
$ws.Range("H4").Value = 1.418309090943296
$ws.Range("I4").Value = 1.0470802714242
$ws.Range("J4").Value = 0.08451451135559829
$ws.Range("M4").Value = 0.5800121765554493
$ws.Range("N4").Value = 1.663908204476719
$ws.Range("B5").Value = 1.083125769798642
$ws.Range("C5").Value = 0.3785569010127006
$ws.Range("D5").Value = 0.5600372932010345
$ws.Range("E5").Value = 0.1997871833376053
$ws.Range("G5").Value = 1.604293392242766
$ws.Range("H5").Value = 1.418295156328455
$ws.Range("I5").Value = 1.04758792969394
$ws.Range("J5").Value = 0.08462726809698573
$ws.Range("M5").Value = 0.5737230017987756
$ws.Range("N5").Value = 1.669517629589464
$ws.Range("B6").Value = 1.079316690229803
$ws.Range("C6").Value = 0.3772354153699666
$ws.Range("D6").Value = 0.5599317761273142
$ws.Range("E6").Value = 0.1997876851751172
$ws.Range("G6").Value = 1.603644791108408
$ws.Range("H6").Value = 1.4183046212639
$ws.Range("I6").Value = 1.047681811702184
$ws.Range("J6").Value = 0.08464644169858637
$ws.Range("M6").Value = 0.5726837641484863
$ws.Range("N6").Value = 1.670458595937646
$ws.Range("B7").Value = 1.10580514425186
$ws.Range("C7").Value = 0.3864270875540115
$ws.Range("D7").Value = 0.5606883285963988
$ws.Range("E7").Value = 0.1997920586364827
$ws.Range("G7").Value = 1.608246848794181
$ws.Range("H7").Value = 1.418308113160492
$ws.Range("I7").Value = 1.047086474989051
$ws.Range("J7").Value = 0.08451600183668972
$ws.Range("M7").Value = 0.5799270182419534
$ws.Range("N7").Value = 1.663983216437648
$ws.Range("B8").Value = 1.223913654613909
$ws.Range("C8").Value = 0.4274632545516965
$ws.Range("D8").Value = 0.5646246647882691
$ws.Range("E8").Value = 0.2000058353068255
$ws.Range("G8").Value = 1.631038737500631
$ws.Range("H8").Value = 1.42003827838721
$ws.Range("I8").Value = 1.045855892339432
$ws.Range("J8").Value = 0.08400587162804563
$ws.Range("M8").Value = 0.6126285320867808
$ws.Range("N8").Value = 1.636794938189851
$ws.Range("B9").Value = 1.459937623822952
$ws.Range("C9").Value = 0.509669685426843
$ws.Range("D9").Value = 0.5744963498230078
$ws.Range("E9").Value = 0.201124908897782
$ws.Range("G9").Value = 1.684744257804539
$ws.Range("H9").Value = 1.429622226478273
$ws.Range("I9").Value = 1.048494153225448
$ws.Range("J9").Value = 0.08324086570694433
$ws.Range("M9").Value = 0.6794244396800764
$ws.Range("N9").Value = 1.588475108345428
$ws.Range("B10").Value = 1.636085331372556
$ws.Range("C10").Value = 0.5711575545572032
$ws.Range("D10").Value = 0.583067129325741
$ws.Range("E10").Value = 0.2023750873700187
$ws.Range("G10").Value = 1.72978109450537
$ws.Range("H10").Value = 1.440468075995938
$ws.Range("I10").Value = 1.053544592885814
$ws.Range("J10").Value = 0.08282240298272825
$ws.Range("M10").Value = 0.7301469061413144
$ws.Range("N10").Value = 1.556031936483002
$ws.Range("B11").Value = 1.716826953594307
$ws.Range("C11").Value = 0.5993745305219704
$ws.Range("D11").Value = 0.587253807435502
$ws.Range("E11").Value = 0.203037090495183
$ws.Range("G11").Value = 1.751500230349876
$ws.Range("H11").Value = 1.44623509589789
$ws.Range("I11").Value = 1.056525183187347
$ws.Range("J11").Value = 0.08266324316813822
$ws.Range("M11").Value = 0.753584577605622
$ws.Range("N11").Value = 1.541939119867788
$ws.Range("B12").Value = 1.747490190891199
$ws.Range("C12").Value = 0.6100954077349456
$ws.Range("D12").Value = 0.5888806643713167
$ws.Range("E12").Value = 0.2033012102760772
$ws.Range("G12").Value = 1.75990325778065
$ws.Range("H12").Value = 1.448539241941347
$ws.Range("I12").Value = 1.057752652241533
$ws.Range("J12").Value = 0.0826074626524651
$ws.Range("M12").Value = 0.7625124083713501
$ws.Range("N12").Value = 1.536698535598706
$ws.Range("B13").Value = 1.740882389900037
$ws.Range("C13").Value = 0.6077848803355437
$ws.Range("D13").Value = 0.588528446862199
$ws.Range("E13").Value = 0.2032437295906533
$ws.Range("G13").Value = 1.758085556637383
$ws.Range("H13").Value = 1.448037644572906
$ws.Range("I13").Value = 1.057483892497885
$ws.Range("J13").Value = 0.08261927624424459
$ws.Range("M13").Value = 0.7605873026174379
$ws.Range("N13").Value = 1.537822911641687
$ws.Range("B14").Value = 1.719347868648128
$ws.Range("C14").Value = 0.6002558254954238
$ws.Range("D14").Value = 0.5873868187503888
$ws.Range("E14").Value = 0.203058550450919
$ws.Range("G14").Value = 1.75218796803253
$ws.Range("H14").Value = 1.446422245542578
$ws.Range("I14").Value = 1.056624184836537
$ws.Range("J14").Value = 0.08265856406896432
$ws.Range("M14").Value = 0.7543180224930666
$ws.Range("N14").Value = 1.541506046777599
$ws.Range("B15").Value = 1.706168838566782
$ws.Range("C15").Value = 0.595648725033584
$ws.Range("D15").Value = 0.586692939055439
$ws.Range("E15").Value = 0.2029468728909976
$ws.Range("G15").Value = 1.748598806124988
$ws.Range("H15").Value = 1.44544844836193
$ws.Range("I15").Value = 1.056110470190362
$ws.Range("J15").Value = 0.08268321380674237
$ws.Range("M15").Value = 0.7504847491508002
$ws.Range("N15").Value = 1.543774592980643
$ws.Range("B16").Value = 1.630820995670547
$ws.Range("C16").Value = 0.5693184870953587
$ws.Range("D16").Value = 0.5827993169422996
$ws.Range("E16").Value = 0.2023337027156309
$ws.Range("G16").Value = 1.728386597127013
$ws.Range("H16").Value = 1.440107995203732
$ws.Range("I16").Value = 1.053363597492542
$ws.Range("J16").Value = 0.08283343256148967
$ws.Range("M16").Value = 0.7286225346758783
$ws.Range("N16").Value = 1.556966359409625
$ws.Range("B17").Value = 1.584754239253584
$ws.Range("C17").Value = 0.5532289640132717
$ws.Range("D17").Value = 0.5804844635449911
$ws.Range("E17").Value = 0.2019814495256611
$ws.Range("G17").Value = 1.716303474620048
$ws.Range("H17").Value = 1.437045545274572
$ws.Range("I17").Value = 1.051853834470755
$ws.Range("J17").Value = 0.08293358041251153
$ws.Range("M17").Value = 0.715304072295865
$ws.Range("N17").Value = 1.565229823980019
$ws.Range("B18").Value = 1.558315366928468
$ws.Range("C18").Value = 0.5439978125033349
$ws.Range("D18").Value = 0.5791801046986791
$ws.Range("E18").Value = 0.2017876226104818
$ws.Range("G18").Value = 1.709469429445875
$ws.Range("H18").Value = 1.435362498994863
$ws.Range("I18").Value = 1.051049723801
$ws.Range("J18").Value = 0.08299411938680379
$ws.Range("M18").Value = 0.7076778730290982
$ws.Range("N18").Value = 1.570045393394707
$ws.Range("B19").Value = 1.549373488560946
$ws.Range("C19").Value = 0.5408762623815164
$ws.Range("D19").Value = 0.5787431207752718
$ws.Range("E19").Value = 0.2017235035373339
$ws.Range("G19").Value = 1.707175395099767
$ws.Range("H19").Value = 1.4348060978879
$ws.Range("I19").Value = 1.050788486599792
$ws.Range("J19").Value = 0.0830151211098844
$ws.Range("M19").Value = 0.7051016462025075
$ws.Range("N19").Value = 1.571686611640231
$ws.Range("B20").Value = 1.589652170521674
$ws.Range("C20").Value = 0.5549393260711213
$ws.Range("D20").Value = 0.5807280800998456
$ws.Range("E20").Value = 0.2020180386725841
$ws.Range("G20").Value = 1.71757774416497
$ws.Range("H20").Value = 1.437363431338213
$ws.Range("I20").Value = 1.052007896171617
$ws.Range("J20").Value = 0.08292261555634539
$ws.Range("M20").Value = 0.7167183015733229
$ws.Range("N20").Value = 1.564343678431449
$ws.Range("B21").Value = 1.725670685392345
$ws.Range("C21").Value = 0.6024663189674584
$ws.Range("D21").Value = 0.5877210170159515
$ws.Range("E21").Value = 0.2031125773183682
$ws.Range("G21").Value = 1.753915378507969
$ws.Range("H21").Value = 1.44689345852683
$ws.Range("I21").Value = 1.056874016407015
$ws.Range("J21").Value = 0.08264690240262595
$ws.Range("M21").Value = 0.7561580358810289
$ws.Range("N21").Value = 1.540421611044451
$ws.Range("B22").Value = 1.815080511797873
$ws.Range("C22").Value = 0.6337363420778388
$ws.Range("D22").Value = 0.5925329584800068
$ws.Range("E22").Value = 0.2039062318926028
$ws.Range("G22").Value = 1.778704974670603
$ws.Range("H22").Value = 1.453823320072871
$ws.Range("I22").Value = 1.060630358113258
$ws.Range("J22").Value = 0.08249287986471288
$ws.Range("M22").Value = 0.7822402436179772
$ws.Range("N22").Value = 1.52534722189305
$ws.Range("B23").Value = 1.767313696747863
$ws.Range("C23").Value = 0.6170277505344757
$ws.Range("D23").Value = 0.5899425998167374
$ws.Range("E23").Value = 0.2034754717777822
$ws.Range("G23").Value = 1.765378612085016
$ws.Range("H23").Value = 1.450060377120906
$ws.Range("I23").Value = 1.058572634961621
$ws.Range("J23").Value = 0.08257268873918377
$ws.Range("M23").Value = 0.7682916111385794
$ws.Range("N23").Value = 1.533341351517336
$ws.Range("B24").Value = 1.587437674028308
$ws.Range("C24").Value = 0.5541660124548002
$ws.Range("D24").Value = 0.5806178585775683
$ws.Range("E24").Value = 0.2020014696586721
$ws.Range("G24").Value = 1.717001295947568
$ws.Range("H24").Value = 1.437219473371954
$ws.Range("I24").Value = 1.051938045959083
$ws.Range("J24").Value = 0.08292756354002151
$ws.Range("M24").Value = 0.7160788327126255
$ws.Range("N24").Value = 1.564744103025401
$ws.Range("B25").Value = 1.395609258068077
$ws.Range("C25").Value = 0.4872415472947296
$ws.Range("D25").Value = 0.5715947634750478
$ws.Range("E25").Value = 0.2007471094606714
$ws.Range("G25").Value = 1.669241972667635
$ws.Range("H25").Value = 1.426363615221135
$ws.Range("I25").Value = 1.047236293658841
$ws.Range("J25").Value = 0.08342261323404543
$ws.Range("M25").Value = 0.661066586079599
$ws.Range("N25").Value = 1.601010848311605
